# Generate Report for Handoff
# Updates the localization-status workbook with new handoff data:
#  - File "155daaac-...md" -> "3819ca96-...md" (now "Ready for handoff")
#  - File "21a6145f-...md" -> "ffff2104bfa4-...md"
#  - Removes the now-unused "Latest Target File" / "Latest Handback File" values
#    (columns F/G) from the per-language detail sheets
#  - Refreshes handoff timestamps

$wb = $excel.ActiveWorkbook

$oldUuid1 = "155daaac-08c1-48d0-964f-fe84de3a97db"
$newUuid1 = "3819ca96-5dc3-45a3-bc6a-3ba87bcd3580"
$oldUuid2 = "21a6145f-f2ad-4ad9-ae1b-10e2d89eec22"
$newUuid2 = "ffff2104bfa4-17a4-4de5-833c-a10edac9086d"
$newHash  = "f7e8e85763eb69c36778c29e41793ccfa9dddb62"

$newMd1  = "$newUuid1.md"
$newMd2  = "$newUuid2.md"
$newZhXlf = "$newUuid1.$newHash.zh-cn.xlf"
$newDeXlf = "$newUuid1.$newHash.de-de.xlf"

$status        = "Ready for handoff"
$overviewDate  = "2016-51-19 16:51:56"
$handoffDateZh = "2016-03-19 16:51:53"
$handoffDateDe = "2016-03-19 16:51:56"
$handbackDate  = "0001-01-01 00:00:00"

$urlBase1 = "https://github.com/OpenLocalizationTest/oltest/blob/e8ff0dd60b3ae1acfffdf6aa37c18561b32bc975/e2e"
$urlZh    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59089f56b5048ba24c9bf369d5af4cdfe414d01e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$urlDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e4df191ceb0cabbc617e4fc490843ead14c2481/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$urlMd1 = "$urlBase1/$newMd1"
$urlMd2 = "$urlBase1/$newMd2"
$urlZhXlf = "$urlZh/$newZhXlf"
$urlDeXlf = "$urlDe/$newDeXlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $overviewDate

$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $overviewDate

$ws1.Hyperlinks.Add($ws1.Range("A2"), $urlMd1, "", "", $newMd1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $urlMd2, "", "", $newMd2)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

$ws2.Range("C2").Value = $status
$ws2.Range("E2").Value = $handoffDateZh
$ws2.Range("F2").Clear()
$ws2.Range("G2").Clear()
$ws2.Range("H2").Value = $handbackDate

$ws2.Range("C3").Value = $status
$ws2.Range("E3").Value = $handoffDateZh
$ws2.Range("F3").Clear()
$ws2.Range("G3").Clear()
$ws2.Range("H3").Value = $handbackDate

$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlMd1, "", "", $newMd1)
$ws2.Hyperlinks.Add($ws2.Range("B2"), $urlMd1, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $urlZhXlf, "", "", $newZhXlf)

$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlMd2, "", "", $newMd2)
$ws2.Hyperlinks.Add($ws2.Range("B3"), $urlMd2, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $urlZhXlf, "", "", $newZhXlf)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Range("C2").Value = $status
$ws3.Range("E2").Value = $handoffDateDe
$ws3.Range("F2").Clear()
$ws3.Range("G2").Clear()
$ws3.Range("H2").Value = $handbackDate

$ws3.Range("C3").Value = $status
$ws3.Range("E3").Value = $handoffDateDe
$ws3.Range("F3").Clear()
$ws3.Range("G3").Clear()
$ws3.Range("H3").Value = $handbackDate

$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlMd1, "", "", $newMd1)
$ws3.Hyperlinks.Add($ws3.Range("B2"), $urlMd1, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $urlDeXlf, "", "", $newDeXlf)

$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlMd2, "", "", $newMd2)
$ws3.Hyperlinks.Add($ws3.Range("B3"), $urlMd2, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $urlDeXlf, "", "", $newDeXlf)
